# Update the "想去人数" (want-to-go count) figures in column F for the
# first and second worksheets that carry this duplicated dataset:
#   - "展览"   (Exhibitions)
#   - "全部类型" (All types)
# Row 2: 337 -> 338
# Row 6: 44  -> 45
# Row 9: 330 -> 332

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 338
    $ws.Range("F6").Value = 45
    $ws.Range("F9").Value = 332
}
